# ByCoach.xlsx - "Add files via upload" re-edit
# Updates the "Started" (Yes/No) flags in column C for several players so
# that, within each coach's block of rows, the "Yes" rows are grouped
# together, and moves the sheet's active selection down to C84 (the
# position the author was last working at when the file was saved).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("by Coach")
$ws.Activate()

# --- Column C ("Started") value swaps -------------------------------
$updates = @{
    2  = "Yes"
    4  = "Yes"
    5  = "No"
    7  = "No"
    14 = "Yes"
    15 = "Yes"
    18 = "No"
    20 = "Yes"
    21 = "No"
    24 = "No"
    28 = "Yes"
    29 = "Yes"
    32 = "No"
    34 = "No"
    53 = "Yes"
    55 = "Yes"
    56 = "No"
    58 = "No"
    65 = "Yes"
    67 = "No"
    79 = "Yes"
    83 = "No"
}

foreach ($row in $updates.Keys) {
    $ws.Range("C$row").Value = $updates[$row]
}

# --- View state: scroll the frozen pane down and leave C84 selected --
$excel.ActiveWindow.ScrollRow = 61
$ws.Range("C84").Select()
